# Relabel save/load toolbar buttons - update tracking sheet status for the
# related feature rows (the sst also picks up the relabeled text for a few
# other in-flight tracking items that were finished in the same pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19: "Ability to rename column row result values..." -> relabeled /
# reclassified and marked Complete.
$ws.Range("B19").Value = "Ability to rename row/column/aggregator labels"
$ws.Range("D19").Value = "HIGH"
$ws.Range("F19").Value = "Complete"
$ws.Range("G19").Value = Get-Date -Year 2017 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H19").Value = Get-Date -Year 2017 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0

# Row 20: "Refine presentation of charts, building on #18" approach note
# gets extra detail appended.
$ws.Range("E20").Value = "Need improvements to chart visuals: sort entires respecting numeric columns (but see #11), make bar charts more legible. Legends take up too much space at bottom of chart. Can charts be resized? Do charts respect label renaming? Also, charts do not work with length 0 aggregators."

# Row 31: "Drag and drop sorting for row/cols" -> approach finalized, marked
# Complete.
$ws.Range("E31").Value = "Sorting done through a dedicated UI element."
$ws.Range("F31").Value = "Complete"
$ws.Range("G31").Value = Get-Date -Year 2017 -Month 9 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("H31").Value = Get-Date -Year 2017 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0

# Row 32: "Row headers should appear hierarchical like col headers" -> marked
# Complete.
$ws.Range("F32").Value = "Complete"
$ws.Range("G32").Value = Get-Date -Year 2017 -Month 9 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("H32").Value = Get-Date -Year 2017 -Month 9 -Day 26 -Hour 0 -Minute 0 -Second 0

# Row 34: "Relabel save/load buttons" -> marked Complete.
$ws.Range("F34").Value = "Complete"
$ws.Range("G34").Value = Get-Date -Year 2017 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("H34").Value = Get-Date -Year 2017 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0

# Move the saved cursor/view position to match where the author left off.
$ws.Activate()
$ws.Range("E32").Select()
$excel.ActiveWindow.ScrollRow = 23
